$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-498) holds a date serial value that was bumped by one day
# (2023-09-08 -> 2023-09-09, serial 45177 -> 45178).
$ws.Range("C2:C498").Value = 45178
